$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.719.22'
$ws.Range('E2').Value = '  -0.20%  '
$ws.Range('D3').Value = '2.326.97'
$ws.Range('E3').Value = '  +4.50%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '271.44'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -0.88%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '95.08'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +8.59%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.626'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  +0.86%  '
$ws.Range('E8').Value = '  -0.06%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.618'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +2.96%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '44.76'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -0.36%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0945'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +3.19%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '8.04'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +5.55%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.104'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -0.28%  '
$ws.Range('D14').Value = '2.673.81'
$ws.Range('E14').Value = '  +4.28%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '15.65'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +4.18%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.855'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +8.33%  '
$ws.Range('D17').Value = '2.329.68'
$ws.Range('E17').Value = '  +3.65%  '
$ws.Range('D18').Value = '43.636.60'
$ws.Range('E18').Value = '  -0.26%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.0000108'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +3.99%  '
$ws.Range('E20').Value = '  +5.89%  '
$ws.Range('E21').Value = '  +2.55%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '237.80'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +2.39%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '2.26'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -3.75%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '9.53'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +10.02%  '
$ws.Range('E25').Value = '  -0.06%  '
$ws.Range('E26').Value = '  -0.57%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '11.27'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +4.41%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '3.41'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -2.44%  '
$ws.Range('E29').Value = '  -0.38%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '38.60'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -1.10%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '22.50'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +8.58%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '172.89'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +0.15%  '
$ws.Range('E33').Value = '  -0.99%  '
$ws.Range('E34').Value = '  +2.54%  '
$ws.Range('E35').Value = '  +3.11%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.0359'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +1.85%  '
$ws.Range('E37').Value = '  -3.84%  '
$ws.Range('E38').Value = '  +1.90%  '
$ws.Range('E39').Value = '  -2.37%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '2.36'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +9.35%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.234'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +12.05%  '
$ws.Range('E42').Value = '  +20.74%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '12.03'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -2.52%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '9.09'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +7.44%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '61.72'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -3.07%  '
$ws.Range('E46').Value = '  -0.31%  '
$ws.Range('E47').Value = '  +5.25%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '100.60'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +0.48%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.21'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +2.09%  '
$ws.Range('D50').Value = '2.553.62'
$ws.Range('E50').Value = '  +4.25%  '
$ws.Range('B51').Value = 'Maker'
$ws.Range('C51').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D51').Value = '1.366.55'
$ws.Range('E51').Value = '  +4.94%  '
